$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column W, shifting the existing "2 Marks" column to X
$ws.Columns("W:W").Insert()

# Set the new column header
$ws.Range("W1").Value = "Total"

# Fill W2:W50 with 0 (Total column values)
$ws.Range("W2:W50").Value = 0
